$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fill in new data columns (D through N)
$ws.Range("D2").Value = "citizenship"
$ws.Range("E2").Value = "citizenship"
$ws.Range("F2").Value = 24
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("L2").Value = 10
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 10

# Update the active selection / view state to match
# (topLeftCell moved to G1, selection moved to P2)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("P2").Select()
